$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F20").Value = "sadsad"
$ws.Range("E20").Value = "asdasdsad"
$ws.Range("E26").Value = "asd"

$ws.Range("E26").Select() | Out-Null
